$wb = $excel.ActiveWorkbook

# --- "Producto" sheet: remove the "Precio" column (column C) ---
$wsProducto = $wb.Worksheets.Item("Producto")
$wsProducto.Range("C1").EntireColumn.Delete()

# Selection on "Producto" sheet becomes a range selection A1:C4 (no explicit active cell)
$wsProducto.Range("A1:C4").Select()

# --- "ProductoPorFabricante" sheet: move selection from E9 to F9 ---
$wsPPF = $wb.Worksheets.Item("ProductoPorFabricante")
$wsPPF.Activate()
$wsPPF.Range("F9").Select()
